{"js": "const oldText = \"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03c4\u03bf\u03c5 \u039a\u03cd\u03ba\u03bd\u03bf\u03c5 2022: 10-19 \u0391\u03c5\u03b3\u03bf\u03cd\u03c3\u03c4\u03bf\u03c5, 9-18 \u03a3\u03b5\u03c0\u03c4\u03b5\u03bc\u03b2\u03c1\u03af\u03bf\u03c5, 8-17 \u039f\u03ba\u03c4\u03c9\u03b2\u03c1\u03af\u03bf\u03c5\";\nconst newText = \"2022 \u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03c4\u03bf\u03c5 \u039a\u03cd\u03ba\u03bd\u03bf\u03c5: 10-19 \u0391\u03c5\u03b3\u03bf\u03cd\u03c3\u03c4\u03bf\u03c5, 9-18 \u03a3\u03b5\u03c0\u03c4\u03b5\u03bc\u03b2\u03c1\u03af\u03bf\u03c5, 8-17 \u039f\u03ba\u03c4\u03c9\u03b2\u03c1\u03af\u03bf\u03c5\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03c4\u03bf\u03c5 \u039a\u03cd\u03ba\u03bd\u03bf\u03c5 2022: 10-19 \u0391\u03c5\u03b3\u03bf\u03cd\u03c3\u03c4\u03bf\u03c5, 9-18 \u03a3\u03b5\u03c0\u03c4\u03b5\u03bc\u03b2\u03c1\u03af\u03bf\u03c5, 8-17 \u039f\u03ba\u03c4\u03c9\u03b2\u03c1\u03af\u03bf\u03c5\"\n$newText = \"2022 \u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03c4\u03bf\u03c5 \u039a\u03cd\u03ba\u03bd\u03bf\u03c5: 10-19 \u0391\u03c5\u03b3\u03bf\u03cd\u03c3\u03c4\u03bf\u03c5, 9-18 \u03a3\u03b5\u03c0\u03c4\u03b5\u03bc\u03b2\u03c1\u03af\u03bf\u03c5, 8-17 \u039f\u03ba\u03c4\u03c9\u03b2\u03c1\u03af\u03bf\u03c5\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newText\n$find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n"}
